$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet and add the two new song-book sheets ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Only Believe Song Book"

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Icilongo"

$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Difelo"

# --- Selection / active sheet bookkeeping ---
# First sheet's selection moves off the old A2:D17 block onto D2.
$ws1.Range("D2").Select() | Out-Null
# Last sheet (Difelo) becomes the active / selected tab.
$ws3.Select() | Out-Null

# --- Normal cell style: point its font back at the plain black font record ---
$normalStyle = $wb.Styles.Item(1)
$normalStyle.Font.ColorIndex = 1
